# Update gh-pages to output generated at 456a3b4
#
# 1) "展览" (sheet1): the "南宁·第一届ANE·DACG动漫嘉年华（取消）" event's attendee
#    count (column F) bumps from 1049 to 1048, and the other four events each get
#    a small attendee-count bump too.
# 2) "演出" (sheet2): the cancelled concert row (row 2, the only data row) is
#    removed entirely, leaving just the header row.
# 3) "全部类型" (sheet4): the same cancelled-concert row is removed (shifting the
#    remaining five events up one row), and the attendee counts get the same
#    bump as in "展览".
#
# Row shifts use Range.Copy (cell-to-cell), NOT Range.Value assignment, so that
# the "开始时间" text cells (e.g. "2024-04-11") are carried over as literal text
# instead of being re-parsed/coerced into date serials by Excel's auto-detect
# (which would also silently mint a new number-format style on the cell).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) 展览 — bump "想去人数" (column F) for rows 2..6
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1048
$wsExpo.Range("F3").Value = 424
$wsExpo.Range("F4").Value = 3148
$wsExpo.Range("F5").Value = 78
$wsExpo.Range("F6").Value = 637

# ---------------------------------------------------------------------------
# 2) 演出 — drop the cancelled concert (only data row, row 2)
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Rows("2:2").Delete()

# ---------------------------------------------------------------------------
# 3) 全部类型 — drop the cancelled concert row (row 2) by shifting the five
#    remaining events up a row (B:I only — column A is a plain 0-based
#    sequence number and must stay 1,2,3,4,5), then delete the now-duplicated
#    last row, then apply the same attendee-count bump as 展览.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("B3:I3").Copy($wsAll.Range("B2:I2"))
$wsAll.Range("B4:I4").Copy($wsAll.Range("B3:I3"))
$wsAll.Range("B5:I5").Copy($wsAll.Range("B4:I4"))
$wsAll.Range("B6:I6").Copy($wsAll.Range("B5:I5"))
$wsAll.Range("B7:I7").Copy($wsAll.Range("B6:I6"))
$wsAll.Rows("7:7").Delete()

$wsAll.Range("F2").Value = 1048
$wsAll.Range("F3").Value = 424
$wsAll.Range("F4").Value = 3148
$wsAll.Range("F5").Value = 78
$wsAll.Range("F6").Value = 637
